# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1680
$ws1.Range("F4").Value = 774
$ws1.Range("F5").Value = 1116
$ws1.Range("F7").Value = 11810
$ws1.Range("F8").Value = 37
$ws1.Range("F10").Value = 470
$ws1.Range("F13").Value = 840
$ws1.Range("F14").Value = 13436
$ws1.Range("F15").Value = 13324

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1680
$ws4.Range("F4").Value = 774
$ws4.Range("F5").Value = 1116
$ws4.Range("F7").Value = 11811
$ws4.Range("F8").Value = 37
$ws4.Range("F10").Value = 470
$ws4.Range("F13").Value = 840
$ws4.Range("F14").Value = 13436
$ws4.Range("F15").Value = 13324
